# Update Name of Algo
# Applies updated imputed values for the KNN result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value  = -12.934
$ws.Range("E6").Value  = 12.857
$ws.Range("E7").Value  = 13.051
$ws.Range("C8").Value  = -12.666
$ws.Range("E8").Value  = 13.084
$ws.Range("A12").Value = -21.882
$ws.Range("C12").Value = -13.002
$ws.Range("C14").Value = -11.675
$ws.Range("E19").Value = 12.614
$ws.Range("E21").Value = 12.93
$ws.Range("C22").Value = -12.929
$ws.Range("E24").Value = 12.81
